{"js": "const replacements = [\n  [\"599\u00d74=\", \"935\u00d79=\"],\n  [\"562\u00d76=\", \"492\u00d72=\"],\n  [\"506\u00d78=\", \"258\u00d72=\"],\n  [\"588\u00d79=\", \"381\u00d76=\"],\n  [\"135\u00d73=\", \"752\u00d79=\"],\n  [\"268\u00d74=\", \"624\u00d78=\"],\n  [\"292\u00d76=\", \"516\u00d72=\"],\n  [\"933\u00d76=\", \"582\u00d74=\"],\n  [\"240\u00d79=\", \"300\u00d76=\"],\n  [\"901\u00d72=\", \"987\u00d77=\"],\n  [\"593\u00d76=\", \"145\u00d76=\"],\n  [\"510\u00d76=\", \"787\u00d78=\"],\n  [\"935\u00d74=\", \"522\u00d76=\"],\n  [\"684\u00d72=\", \"671\u00d77=\"],\n  [\"907\u00d72=\", \"647\u00d72=\"],\n  [\"853\u00d79=\", \"998\u00d73=\"],\n  [\"514\u00d79=\", \"680\u00d76=\"],\n  [\"545\u00d79=\", \"265\u00d75=\"],\n  [\"243\u00d74=\", \"302\u00d73=\"],\n  [\"941\u00d77=\", \"972\u00d77=\"],\n  [\"745\u00d73=\", \"954\u00d75=\"],\n  [\"444\u00d74=\", \"878\u00d79=\"],\n  [\"606\u00d72=\", \"938\u00d76=\"],\n  [\"829\u00d76=\", \"907\u00d74=\"],\n  [\"801\u00d76=\", \"709\u00d75=\"],\n];\n\nfor (const [oldText, newText] of replacements) {\n  const results = context.document.body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n  for (const range of results.items) {\n    range.insertText(newText, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n", "ps1": "$d = $word.ActiveDocument\n\n$replacements = @(\n    @(\"599\u00d74=\", \"935\u00d79=\"),\n    @(\"562\u00d76=\", \"492\u00d72=\"),\n    @(\"506\u00d78=\", \"258\u00d72=\"),\n    @(\"588\u00d79=\", \"381\u00d76=\"),\n    @(\"135\u00d73=\", \"752\u00d79=\"),\n    @(\"268\u00d74=\", \"624\u00d78=\"),\n    @(\"292\u00d76=\", \"516\u00d72=\"),\n    @(\"933\u00d76=\", \"582\u00d74=\"),\n    @(\"240\u00d79=\", \"300\u00d76=\"),\n    @(\"901\u00d72=\", \"987\u00d77=\"),\n    @(\"593\u00d76=\", \"145\u00d76=\"),\n    @(\"510\u00d76=\", \"787\u00d78=\"),\n    @(\"935\u00d74=\", \"522\u00d76=\"),\n    @(\"684\u00d72=\", \"671\u00d77=\"),\n    @(\"907\u00d72=\", \"647\u00d72=\"),\n    @(\"853\u00d79=\", \"998\u00d73=\"),\n    @(\"514\u00d79=\", \"680\u00d76=\"),\n    @(\"545\u00d79=\", \"265\u00d75=\"),\n    @(\"243\u00d74=\", \"302\u00d73=\"),\n    @(\"941\u00d77=\", \"972\u00d77=\"),\n    @(\"745\u00d73=\", \"954\u00d75=\"),\n    @(\"444\u00d74=\", \"878\u00d79=\"),\n    @(\"606\u00d72=\", \"938\u00d76=\"),\n    @(\"829\u00d76=\", \"907\u00d74=\"),\n    @(\"801\u00d76=\", \"709\u00d75=\")\n)\n\nforeach ($pair in $replacements) {\n    $find = $pair[0]\n    $replace = $pair[1]\n    $range = $d.Content\n    $range.Find.ClearFormatting()\n    $range.Find.Replacement.ClearFormatting()\n    $range.Find.Text = $find\n    $range.Find.Replacement.Text = $replace\n    $range.Find.Forward = $true\n    $range.Find.Wrap = 1\n    $range.Find.MatchCase = $true\n    $range.Find.MatchWholeWord = $false\n    $range.Find.MatchWildcards = $false\n    $range.Find.Execute(\n        $find, $true, $false, $false, $false, $false, $true, 1, $false, $replace, 2\n    )\n}\n"}
